# "finish the last step" -- add the final hlf15 benchmark rows to the
# "15" sheet, fix two mislabeled qft15 entries, and leave the "12"
# sheet's selection parked on the newly-finished M16:U19 block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "15": fix two mislabeled A-column entries (they were stamped
# with the "qft12(...)" label instead of "qft15(...)").
# ---------------------------------------------------------------------
$ws15 = $wb.Worksheets.Item("15")

$ws15.Range("A9").Value = "qft15(38)"
$ws15.Range("A10").Value = "qft15(33)"

# ---------------------------------------------------------------------
# Sheet "15": append the last four benchmark rows (hlf15 family),
# mirroring the formula pattern already used in B8:I15.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row = 16; Label = "hlf15(41)"; B = 52; C = 46; E = 35; F = 37 },
    @{ Row = 17; Label = "hlf15(34)"; B = 52; C = 46; E = 37; F = 39 },
    @{ Row = 18; Label = "hlf15(33)"; B = 52; C = 46; E = 37; F = 40 },
    @{ Row = 19; Label = "hlf15(29)"; B = 52; C = 46; E = 46; F = 49 }
)

foreach ($r in $newRows) {
    $i = $r.Row
    $ws15.Range("A$i").Value = $r.Label
    $ws15.Range("B$i").Value = $r.B
    $ws15.Range("C$i").Value = $r.C
    $ws15.Range("D$i").Formula = "=B$i-C$i"
    $ws15.Range("E$i").Value = $r.E
    $ws15.Range("F$i").Value = $r.F
    $ws15.Range("G$i").Formula = "=(C$i-E$i)/B$i"
    $ws15.Range("H$i").Formula = "=(D$i-F$i+E$i)/B$i"
    $ws15.Range("I$i").Formula = "= 1 -F$i/B$i"
}

# Move the selection off the data, matching the end-of-entry state.
$ws15.Activate()
$ws15.Range("I25").Select()

# ---------------------------------------------------------------------
# Sheet "12": the last thing the author did was select the block that
# now mirrors the finished sheet ("15"), M16:U19.
# ---------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("12")
$ws12.Activate()
$ws12.Range("M16:U19").Select()

# Leave "15" as the active tab (it was tabSelected in the original file).
$ws15.Activate()
